$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as text, even if it looks numeric (e.g. "2020"),
# so it is stored as a string rather than being auto-converted to a number.
function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.ClearFormats()
}

# Remove rows 7-11 entirely (dimension shrinks to A1:H6)
$ws.Rows("7:11").Delete()

# Row 2
$ws.Range("B2").Value = "supervisor"
$ws.Range("C2").Value = "TRUONG LONG EXPORT -"
$ws.Range("D2").Value = ""
Set-TextValue $ws.Range("E2") "2020"
Set-TextValue $ws.Range("F2") "2020"

# Row 3
$ws.Range("B3").Value = "supervisor"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
Set-TextValue $ws.Range("E3") "2018"
Set-TextValue $ws.Range("F3") "2020"

# Row 4
$ws.Range("B4").Value = "supervisor"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
Set-TextValue $ws.Range("E4") "2014"
Set-TextValue $ws.Range("F4") "2018"

# Row 5
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
Set-TextValue $ws.Range("E5") "2014"
Set-TextValue $ws.Range("F5") "2014"

# Row 6
$ws.Range("C6").Value = ""
Set-TextValue $ws.Range("E6") "2012"
Set-TextValue $ws.Range("F6") "2013"
